# Apply the crypto price/volume refresh described in the commit.
# Column D holds prices as text (e.g. "29.832.78", "1.000") and column E
# holds percentage-change text; both must stay plain text, matching the
# workbook's original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.832.78'
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.887.26'
$ws.Range("E3").Value = '  -0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7486'
$ws.Range("E5").Value = '  -4.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.24'
$ws.Range("E6").Value = '  -0.65%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3129'
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.20'
$ws.Range("E9").Value = '  -2.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07126'
$ws.Range("E10").Value = '  -2.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08514'
$ws.Range("E11").Value = '  +5.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7596'
$ws.Range("E12").Value = '  -1.88%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.362'
$ws.Range("E13").Value = '  -2.06%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.36'
$ws.Range("E14").Value = '  -0.56%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.811.75'
$ws.Range("E15").Value = '  -1.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.161'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.829.36'
$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.69'
$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.40'
$ws.Range("E19").Value = '  -1.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007795'
$ws.Range("E20").Value = '  -0.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.146.01'
$ws.Range("E21").Value = '  +1.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.993'
$ws.Range("E23").Value = '  -1.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1586'
$ws.Range("E25").Value = '  -0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.357'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.26'
$ws.Range("E27").Value = '  -1.14%  '

$ws.Range("E28").Value = '  +0.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.027'
$ws.Range("E29").Value = '  +0.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.501'
$ws.Range("E30").Value = '  +4.43%  '

$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.511'
$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.121'
$ws.Range("E33").Value = '  +1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05412'
$ws.Range("E34").Value = '  -2.66%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.239'
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7467'
$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.710'
$ws.Range("E38").Value = '  +0.92%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01942'
$ws.Range("E39").Value = '  +0.21%  '

$ws.Range("E40").Value = '  -1.01%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4452'
$ws.Range("E41").Value = '  -0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.084'
$ws.Range("E42").Value = '  +2.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.090.33'
$ws.Range("E43").Value = '  -4.33%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.37'
$ws.Range("E44").Value = '  -2.14%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8550'
$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.32'
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.693'
$ws.Range("E48").Value = '  +2.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.859'
$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.038'
$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.034.94'
$ws.Range("E51").Value = '  -0.04%  '
